$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 2936.4285
$ws.Range("I2").Value2 = 5196.25
$ws.Range("K2").Value2 = 5196.25
$ws.Range("M2").Value2 = -5083.25
$ws.Range("H4").Value2 = 1442.8572
$ws.Range("I4").Value2 = 1350.1666
$ws.Range("J4").Value2 = 1999
$ws.Range("K4").Value2 = 1350.1666
$ws.Range("L4").Value2 = 1999
$ws.Range("M4").Value2 = -1236.1666
$ws.Range("N4").Value2 = -2227
$ws.Range("H70").Value2 = 8244.294
$ws.Range("I70").Value2 = 7595.8335
$ws.Range("J70").Value2 = 9800.6
$ws.Range("K70").Value2 = 22787.5005
$ws.Range("L70").Value2 = 29401.8
$ws.Range("M70").Value2 = -22517.5005
$ws.Range("N70").Value2 = -29941.8
$ws.Range("H73").Value2 = 8244.294
$ws.Range("I73").Value2 = 7595.8335
$ws.Range("J73").Value2 = 9800.6
$ws.Range("K73").Value2 = 22787.5005
$ws.Range("L73").Value2 = 29401.8
$ws.Range("M73").Value2 = -21851.5005
$ws.Range("N73").Value2 = -31273.8
$ws.Range("H86").Value2 = 107799.9
$ws.Range("I86").Value2 = 30000
$ws.Range("J86").Value2 = 116444.336
$ws.Range("K86").Value2 = 30000
$ws.Range("L86").Value2 = 116444.336
$ws.Range("M86").Value2 = -28877
$ws.Range("N86").Value2 = -118690.336
$ws.Range("H89").Value2 = 107799.9
$ws.Range("I89").Value2 = 30000
$ws.Range("J89").Value2 = 116444.336
$ws.Range("K89").Value2 = 150000
$ws.Range("L89").Value2 = 582221.6799999999
$ws.Range("M89").Value2 = -144384
$ws.Range("N89").Value2 = -593453.6799999999
$ws.Range("H137").Value2 = 4768.125
$ws.Range("I137").Value2 = 10425.5
$ws.Range("K137").Value2 = 31276.5
$ws.Range("M137").Value2 = -28726.5
$ws.Range("H138").Value2 = 3812.3333
$ws.Range("J138").Value2 = 5689.4287
$ws.Range("L138").Value2 = 17068.2861
$ws.Range("N138").Value2 = -27348.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value2 = 33334564
$ws.Range("J2").Value2 = 1050
$ws.Range("K2").Value2 = 33334564
$ws.Range("L2").Value2 = 1050
$ws.Range("M2").Value2 = -33334451
$ws.Range("N2").Value2 = -1276
$ws.Range("H45").Value2 = 1655.25
$ws.Range("I45").Value2 = 1655.25
$ws.Range("K45").Value2 = 1655.25
$ws.Range("M45").Value2 = -1278.25
$ws.Range("H61").Value2 = 9332.166999999999
$ws.Range("I61").Value2 = 9332.166999999999
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 9332.166999999999
$ws.Range("L61").Value2 = 0
$ws.Range("M61").Value2 = -9120.166999999999
$ws.Range("N61").ClearContents() | Out-Null
$ws.Range("H98").Value2 = 28089
$ws.Range("J98").Value2 = 28089
$ws.Range("L98").Value2 = 28089
$ws.Range("N98").Value2 = -34079
$ws.Range("H110").Value2 = 3832836.5
$ws.Range("J110").Value2 = 1959.8667
$ws.Range("L110").Value2 = 1959.8667
$ws.Range("N110").Value2 = -6049.8667
$ws.Range("I116").Value2 = 33334564
$ws.Range("J116").Value2 = 1050
$ws.Range("K116").Value2 = 33334564
$ws.Range("L116").Value2 = 1050
$ws.Range("M116").Value2 = -33332270
$ws.Range("N116").Value2 = -5638
$ws.Range("H122").Value2 = 3805.3914
$ws.Range("I122").Value2 = 3121.2307
$ws.Range("K122").Value2 = 9363.6921
$ws.Range("M122").Value2 = -6913.6921
$ws.Range("H136").Value2 = 9332.166999999999
$ws.Range("I136").Value2 = 9332.166999999999
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 27996.501
$ws.Range("L136").Value2 = 0
$ws.Range("M136").Value2 = -25446.501
$ws.Range("N136").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value2 = 33334564
$ws.Range("J3").Value2 = 1050
$ws.Range("K3").Value2 = 33334564
$ws.Range("L3").Value2 = 1050
$ws.Range("M3").Value2 = -33334450
$ws.Range("N3").Value2 = -1278
$ws.Range("H35").Value2 = 15166.5
$ws.Range("J35").Value2 = 15166.5
$ws.Range("L35").Value2 = 15166.5
$ws.Range("N35").Value2 = -15786.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 75798
$ws.Range("I22").Value2 = 105017.3
$ws.Range("J22").Value2 = 2749.75
$ws.Range("K22").Value2 = 105017.3
$ws.Range("L22").Value2 = 2749.75
$ws.Range("M22").Value2 = -104667.3
$ws.Range("N22").Value2 = -3449.75
$ws.Range("H31").Value2 = 4498.2856
$ws.Range("I31").Value2 = 4678.8335
$ws.Range("K31").Value2 = 4678.8335
$ws.Range("M31").Value2 = -4383.8335
$ws.Range("H34").Value2 = 4498.2856
$ws.Range("I34").Value2 = 4678.8335
$ws.Range("K34").Value2 = 4678.8335
$ws.Range("M34").Value2 = -4476.8335
$ws.Range("H43").Value2 = 18828.5
$ws.Range("J43").Value2 = 18828.5
$ws.Range("L43").Value2 = 18828.5
$ws.Range("N43").Value2 = -19196.5
$ws.Range("H58").Value2 = 3010.5173
$ws.Range("I58").Value2 = 2607.2104
$ws.Range("J58").Value2 = 3776.8
$ws.Range("K58").Value2 = 2607.2104
$ws.Range("L58").Value2 = 3776.8
$ws.Range("M58").Value2 = -2404.2104
$ws.Range("N58").Value2 = -4182.8
$ws.Range("H92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("N92").ClearContents() | Out-Null
$ws.Range("H101").Value2 = 18828.5
$ws.Range("J101").Value2 = 18828.5
$ws.Range("L101").Value2 = 18828.5
$ws.Range("N101").Value2 = -25318.5
$ws.Range("H119").Value2 = 0
$ws.Range("J119").Value2 = 0
$ws.Range("L119").Value2 = 0
$ws.Range("N119").ClearContents() | Out-Null
$ws.Range("H130").Value2 = 40779.25
$ws.Range("J130").Value2 = 40779.25
$ws.Range("L130").Value2 = 40779.25
$ws.Range("N130").Value2 = -50819.25
$ws.Range("H134").Value2 = 2674.36
$ws.Range("I134").Value2 = 2557.7144
$ws.Range("K134").Value2 = 7673.1432
$ws.Range("M134").Value2 = -5138.1432
$ws.Range("H136").Value2 = 3010.5173
$ws.Range("I136").Value2 = 2607.2104
$ws.Range("J136").Value2 = 3776.8
$ws.Range("K136").Value2 = 7821.6312
$ws.Range("L136").Value2 = 11330.4
$ws.Range("M136").Value2 = -5271.6312
$ws.Range("N136").Value2 = -16430.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 1991.6
$ws.Range("I68").Value2 = 1510
$ws.Range("J68").Value2 = 2232.4
$ws.Range("K68").Value2 = 4530
$ws.Range("L68").Value2 = 6697.200000000001
$ws.Range("M68").Value2 = -3719
$ws.Range("N68").Value2 = -8319.200000000001
$ws.Range("H71").Value2 = 1991.6
$ws.Range("I71").Value2 = 1510
$ws.Range("J71").Value2 = 2232.4
$ws.Range("K71").Value2 = 13590
$ws.Range("L71").Value2 = 20091.6
$ws.Range("M71").Value2 = -9534
$ws.Range("N71").Value2 = -28203.6
$ws.Range("H113").Value2 = 1051.7142
$ws.Range("I113").Value2 = 854.3333
$ws.Range("J113").Value2 = 1199.75
$ws.Range("K113").Value2 = 2562.9999
$ws.Range("L113").Value2 = 3599.25
$ws.Range("M113").Value2 = -392.9998999999998
$ws.Range("N113").Value2 = -7939.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value2 = 49998
$ws.Range("J123").Value2 = 49998
$ws.Range("L123").Value2 = 49998
$ws.Range("N123").Value2 = -54898
$ws.Range("H126").Value2 = 4501.4
$ws.Range("I126").Value2 = 4002.8
$ws.Range("K126").Value2 = 12008.4
$ws.Range("M126").Value2 = -9538.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1109.6
$ws.Range("J22").Value2 = 999.75
$ws.Range("L22").Value2 = 999.75
$ws.Range("N22").Value2 = -1589.75
$ws.Range("H27").Value2 = 1109.6
$ws.Range("J27").Value2 = 999.75
$ws.Range("L27").Value2 = 999.75
$ws.Range("N27").Value2 = -1213.75
$ws.Range("H51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("N51").ClearContents() | Out-Null
$ws.Range("H93").Value2 = 12112.857
$ws.Range("I93").Value2 = 12965
$ws.Range("K93").Value2 = 12965
$ws.Range("M93").Value2 = -11717
$ws.Range("H97").Value2 = 30344
$ws.Range("J97").Value2 = 30344
$ws.Range("L97").Value2 = 30344
$ws.Range("N97").Value2 = -32326
$ws.Range("H132").Value2 = 2922
$ws.Range("I132").Value2 = 2922
$ws.Range("K132").Value2 = 8766
$ws.Range("M132").Value2 = -6236
$ws.Range("H136").Value2 = 6664
$ws.Range("I136").Value2 = 6664
$ws.Range("K136").Value2 = 19992
$ws.Range("M136").Value2 = -17442

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value2 = 38689
$ws.Range("J47").Value2 = 38689
$ws.Range("L47").Value2 = 38689
$ws.Range("N47").Value2 = -39833
$ws.Range("H95").Value2 = 20000
$ws.Range("J95").Value2 = 20000
$ws.Range("L95").Value2 = 20000
$ws.Range("N95").Value2 = -25492
$ws.Range("H132").Value2 = 1099.1
$ws.Range("I132").Value2 = 956.5714
$ws.Range("J132").Value2 = 1431.6666
$ws.Range("K132").Value2 = 2869.7142
$ws.Range("L132").Value2 = 4294.9998
$ws.Range("M132").Value2 = -339.7142000000003
$ws.Range("N132").Value2 = -9354.9998
